# Add a new user row (Suriya) to the users table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "Suriya"

# Email column becomes a mailto: hyperlink, styled with Excel's built-in
# "Hyperlink" cell style; set the display text explicitly afterwards so the
# visible/stored text is the address itself rather than the full mailto: URL.
$ws.Hyperlinks.Add($ws.Range("B14"), "mailto:suriya@gmail.com") | Out-Null
$ws.Range("B14").Value = "suriya@gmail.com"

$ws.Range("C14").Value = 9941848089
$ws.Range("D14").Value = 19
$ws.Range("E14").Value = "Male"
$ws.Range("F14").Value = "Cheannai"

# Widen the Email/ContactNumber columns so the new data is fully visible.
$ws.Columns("B").ColumnWidth = 22.6
$ws.Columns("C").ColumnWidth = 10.1

# Leave the selection where the editor last left it.
$ws.Range("H7").Select() | Out-Null
